$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.28222890584766
$ws.Range("C2").Value = 0.5481659376659422
$ws.Range("D2").Value = 0.04427475039739415
$ws.Range("E2").Value = 0.0908676865586191
$ws.Range("F2").Value = 3.156006502200995
$ws.Range("I2").Value = 1.56739384145051
$ws.Range("J2").Value = 0.1622475285559801
$ws.Range("M2").Value = 0.5066399124069534
$ws.Range("B3").Value = 1.200853625392824
$ws.Range("C3").Value = 0.5139791758210492
$ws.Range("D3").Value = 0.04430237240150681
$ws.Range("E3").Value = 0.09104512587838087
$ws.Range("F3").Value = 3.109600644876309
$ws.Range("I3").Value = 1.549139059382128
$ws.Range("J3").Value = 0.161989779380562
$ws.Range("M3").Value = 0.4871025092023373
$ws.Range("B4").Value = 1.151708354095263
$ws.Range("C4").Value = 0.493347295189551
$ws.Range("D4").Value = 0.04432452183390856
$ws.Range("E4").Value = 0.09118347383091319
$ws.Range("F4").Value = 3.082774053227766
$ws.Range("I4").Value = 1.538704150164577
$ws.Range("J4").Value = 0.1618953660720308
$ws.Range("M4").Value = 0.4754341351724847
$ws.Range("B5").Value = 1.131886617910197
$ws.Range("C5").Value = 0.4850291901934156
$ws.Range("D5").Value = 0.04433485433507967
$ws.Range("E5").Value = 0.09124725632203301
$ws.Range("F5").Value = 3.072258635838764
$ws.Range("I5").Value = 1.534645015360184
$ws.Range("J5").Value = 0.1618729085125388
$ws.Range("M5").Value = 0.4707613769008319
$ws.Range("B6").Value = 1.128607617262531
$ws.Range("C6").Value = 0.4836533638363392
$ws.Range("D6").Value = 0.04433664899429268
$ws.Range("E6").Value = 0.0912582949277585
$ws.Range("F6").Value = 3.070537650114417
$ws.Range("I6").Value = 1.533982626986948
$ws.Range("J6").Value = 0.1618701457199236
$ws.Range("M6").Value = 0.4699904280438005
$ws.Range("B7").Value = 1.151440200643094
$ws.Range("C7").Value = 0.4932347526296326
$ws.Range("D7").Value = 0.04432465588941081
$ws.Range("E7").Value = 0.09118430402685718
$ws.Range("F7").Value = 3.082630555006389
$ws.Range("I7").Value = 1.538648626900411
$ws.Range("J7").Value = 0.1618949984003066
$ws.Range("M7").Value = 0.4753707841542649
$ws.Range("B8").Value = 1.254000162951002
$ws.Range("C8").Value = 0.5363034081219098
$ws.Range("D8").Value = 0.04428319819876947
$ws.Range("E8").Value = 0.09092277268183935
$ws.Range("F8").Value = 3.139658294903114
$ws.Range("I8").Value = 1.56093820798759
$ws.Range("J8").Value = 0.1621453821331045
$ws.Range("M8").Value = 0.4998352789478773
$ws.Range("B9").Value = 1.461667651512073
$ws.Range("C9").Value = 0.623646004107627
$ws.Range("D9").Value = 0.04424303362808257
$ws.Range("E9").Value = 0.09064272450688726
$ws.Range("F9").Value = 3.264836396595513
$ws.Range("I9").Value = 1.610850791761791
$ws.Range("J9").Value = 0.1631448726829419
$ws.Range("M9").Value = 0.5504225436411758
$ws.Range("B10").Value = 1.618315520526835
$ws.Range("C10").Value = 0.6896348235301275
$ws.Range("D10").Value = 0.04423858127484825
$ws.Range("E10").Value = 0.09057838001977103
$ws.Range("F10").Value = 3.365123124212118
$ws.Range("I10").Value = 1.651397602501618
$ws.Range("J10").Value = 0.1641920994010704
$ws.Range("M10").Value = 0.5892043350370813
$ws.Range("B11").Value = 1.690483181500269
$ws.Range("C11").Value = 0.7200632040289747
$ws.Range("D11").Value = 0.04424200005883705
$ws.Range("E11").Value = 0.09057971415337462
$ws.Range("F11").Value = 3.412592780672696
$ws.Range("I11").Value = 1.670705889722925
$ws.Range("J11").Value = 0.1647371099160964
$ws.Range("M11").Value = 0.6072033333616247
$ws.Range("B12").Value = 1.717943016845368
$ws.Range("C12").Value = 0.7316455331157954
$ws.Range("D12").Value = 0.04424407802898411
$ws.Range("E12").Value = 0.09058461136794094
$ws.Range("F12").Value = 3.430837209423004
$ws.Range("I12").Value = 1.678143188681929
$ws.Range("J12").Value = 0.1649534087570288
$ws.Range("M12").Value = 0.614070774168809
$ws.Range("B13").Value = 1.712023185112514
$ws.Range("C13").Value = 0.7291483997607315
$ws.Range("D13").Value = 0.04424359564936564
$ws.Range("E13").Value = 0.09058336146806134
$ws.Range("F13").Value = 3.426895948897624
$ws.Range("I13").Value = 1.676535819698387
$ws.Range("J13").Value = 0.1649063831865263
$ws.Range("M13").Value = 0.6125894479266094
$ws.Range("B14").Value = 1.692739677315501
$ws.Range("C14").Value = 0.7210148865999031
$ws.Range("D14").Value = 0.04424215531068088
$ws.Range("E14").Value = 0.09058002907764795
$ws.Range("F14").Value = 3.414088359177612
$ws.Range("I14").Value = 1.671315233410567
$ws.Range("J14").Value = 0.1647547059618901
$ws.Range("M14").Value = 0.6077672858467622
$ws.Range("B15").Value = 1.680945127951418
$ws.Range("C15").Value = 0.7160406830379884
$ws.Range("D15").Value = 0.04424137509681758
$ws.Range("E15").Value = 0.09057855960243266
$ws.Range("F15").Value = 3.406278418685076
$ws.Range("I15").Value = 1.668133884873839
$ws.Range("J15").Value = 0.1646630919546013
$ws.Range("M15").Value = 0.6048203020945238
$ws.Range("B16").Value = 1.613617542778741
$ws.Range("C16").Value = 0.6876545737663378
$ws.Range("D16").Value = 0.04423846743398996
$ws.Range("E16").Value = 0.09057890790666789
$ws.Range("F16").Value = 3.362058314913156
$ws.Range("I16").Value = 1.650153262963173
$ws.Range("J16").Value = 0.164157866566363
$ws.Range("M16").Value = 0.5880352648396254
$ws.Range("B17").Value = 1.572547427185555
$ws.Range("C17").Value = 0.6703461511013415
$ws.Range("D17").Value = 0.04423807830434012
$ws.Range("E17").Value = 0.09058695373834702
$ws.Range("F17").Value = 3.33540623705926
$ws.Range("I17").Value = 1.639344934870437
$ws.Range("J17").Value = 0.1638655363571289
$ws.Range("M17").Value = 0.5778298106269659
$ws.Range("B18").Value = 1.549010364526623
$ws.Range("C18").Value = 0.6604293188095198
$ws.Range("D18").Value = 0.04423836686944815
$ws.Range("E18").Value = 0.09059446266563675
$ws.Range("F18").Value = 3.320250504437297
$ws.Range("I18").Value = 1.633209440244983
$ws.Range("J18").Value = 0.1637038518785445
$ws.Range("M18").Value = 0.5719934852050272
$ws.Range("B19").Value = 1.541055754982892
$ws.Range("C19").Value = 0.6570782393871468
$ws.Range("D19").Value = 0.04423855256556664
$ws.Range("E19").Value = 0.09059750018164969
$ws.Range("F19").Value = 3.315148804055866
$ws.Range("I19").Value = 1.631145961996197
$ws.Range("J19").Value = 0.1636502157051822
$ws.Range("M19").Value = 0.5700231651224854
$ws.Range("B20").Value = 1.576910567530604
$ws.Range("C20").Value = 0.6721846700883134
$ws.Range("D20").Value = 0.0442380666957094
$ws.Range("E20").Value = 0.09058579912064069
$ws.Range("F20").Value = 3.338225384131732
$ws.Range("I20").Value = 1.640487088095313
$ws.Range("J20").Value = 0.1638959868817054
$ws.Range("M20").Value = 0.5789127210671055
$ws.Range("B21").Value = 1.698400133904727
$ws.Range("C21").Value = 0.7234022678422889
$ws.Range("D21").Value = 0.04424255710796676
$ws.Range("E21").Value = 0.09058088875307391
$ws.Range("F21").Value = 3.417842943363695
$ws.Range("I21").Value = 1.672845223313487
$ws.Range("J21").Value = 0.1647989877621541
$ws.Range("M21").Value = 0.6091822687481709
$ws.Range("B22").Value = 1.77856760510025
$ws.Range("C22").Value = 0.7572246718008273
$ws.Range("D22").Value = 0.04425005828220119
$ws.Range("E22").Value = 0.09060327562828974
$ws.Range("F22").Value = 3.471445160357888
$ws.Range("I22").Value = 1.694726350248899
$ws.Range("J22").Value = 0.1654469669526222
$ws.Range("M22").Value = 0.6292661339806642
$ws.Range("B23").Value = 1.73571020109847
$ws.Range("C23").Value = 0.739140824568949
$ws.Range("D23").Value = 0.04424563669063675
$ws.Range("E23").Value = 0.09058898814584637
$ws.Range("F23").Value = 3.442692265771655
$ws.Range("I23").Value = 1.682980385206662
$ws.Range("J23").Value = 0.1650958222088192
$ws.Range("M23").Value = 0.6185193630266923
$ws.Range("B24").Value = 1.574937759198178
$ws.Range("C24").Value = 0.671353369782878
$ws.Range("D24").Value = 0.04423807034794081
$ws.Range("E24").Value = 0.09058631214190349
$ws.Range("F24").Value = 3.336950328093394
$ws.Range("I24").Value = 1.639970476672218
$ws.Range("J24").Value = 0.1638822003304696
$ws.Range("M24").Value = 0.5784230409499216
$ws.Range("B25").Value = 1.404778708955973
$ws.Range("C25").Value = 0.5997025805866087
$ws.Range("D25").Value = 0.04424949877493978
$ws.Range("E25").Value = 0.09069361587740232
$ws.Range("F25").Value = 3.229524848490996
$ws.Range("I25").Value = 1.596674322333854
$ws.Range("J25").Value = 0.1628197483952007
$ws.Range("M25").Value = 0.5364553647189538
